# Updated symbol list on Fri Jan  6 22:44:16 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# the coin rows on Sheet1. The source values are stored as plain text
# (e.g. "259.09", "0.80%"), so each cell is explicitly formatted as Text
# before the new value is written - this prevents Excel from silently
# re-interpreting the strings as numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2";  Value = "259.09" }
    @{ Cell = "E2";  Value = "0.80%" }
    @{ Cell = "D3";  Value = "26.96" }
    @{ Cell = "D4";  Value = "4.682" }
    @{ Cell = "E4";  Value = "0.53%" }
    @{ Cell = "D5";  Value = "0.06043" }
    @{ Cell = "E5";  Value = "2.76%" }
    @{ Cell = "D6";  Value = "6.675" }
    @{ Cell = "E6";  Value = "0.74%" }
    @{ Cell = "E7";  Value = "0.02%" }
    @{ Cell = "D8";  Value = "0.9205" }
    @{ Cell = "E8";  Value = "-3.13%" }
    @{ Cell = "E9";  Value = "-0.84%" }
    @{ Cell = "D10"; Value = "0.05283" }
    @{ Cell = "E10"; Value = "28.97%" }
    @{ Cell = "D11"; Value = "0.07085" }
    @{ Cell = "E11"; Value = "-0.04%" }
    @{ Cell = "D12"; Value = "0.03067" }
    @{ Cell = "E12"; Value = "-3.58%" }
    @{ Cell = "D13"; Value = "0.09129" }
    @{ Cell = "E13"; Value = "-0.39%" }
    @{ Cell = "D14"; Value = "0.001536" }
    @{ Cell = "E14"; Value = "-1.13%" }
    @{ Cell = "D15"; Value = "0.0006051" }
    @{ Cell = "E15"; Value = "-94.24%" }
    @{ Cell = "D16"; Value = "0.006066" }
    @{ Cell = "E16"; Value = "-2.18%" }
    @{ Cell = "D17"; Value = "3.469" }
    @{ Cell = "E17"; Value = "-1.72%" }
    @{ Cell = "D18"; Value = "3.174" }
    @{ Cell = "E18"; Value = "-0.95%" }
    @{ Cell = "E19"; Value = "-0.93%" }
    @{ Cell = "E20"; Value = "2.40%" }
    @{ Cell = "E21"; Value = "-0.22%" }
    @{ Cell = "D22"; Value = "4.110" }
    @{ Cell = "E22"; Value = "7.34%" }
    @{ Cell = "D23"; Value = "0.04244" }
    @{ Cell = "E23"; Value = "0.43%" }
    @{ Cell = "D24"; Value = "0.001218" }
    @{ Cell = "E24"; Value = "-0.23%" }
    @{ Cell = "E25"; Value = "-6.29%" }
    @{ Cell = "D26"; Value = "0.0001200" }
    @{ Cell = "E26"; Value = "-0.01%" }
    @{ Cell = "E27"; Value = "-21.36%" }
    @{ Cell = "D40"; Value = "0.03859" }
    @{ Cell = "E40"; Value = "0.76%" }
    @{ Cell = "D41"; Value = "0.1115" }
    @{ Cell = "E41"; Value = "1.33%" }
    @{ Cell = "D43"; Value = "0.01509" }
    @{ Cell = "E43"; Value = "32.20%" }
    @{ Cell = "D44"; Value = "0.002200" }
    @{ Cell = "E44"; Value = "-0.01%" }
    @{ Cell = "D45"; Value = "0.00005185" }
    @{ Cell = "E45"; Value = "-5.26%" }
    @{ Cell = "E46"; Value = "-0.01%" }
    @{ Cell = "D47"; Value = "0.05452" }
    @{ Cell = "E47"; Value = "-22.13%" }
    @{ Cell = "D48"; Value = "0.1321" }
    @{ Cell = "E48"; Value = "-43.26%" }
    @{ Cell = "D49"; Value = "0.00002100" }
    @{ Cell = "E49"; Value = "-0.01%" }
    @{ Cell = "D50"; Value = "0.0002000" }
    @{ Cell = "E50"; Value = "-0.01%" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
